# Commit 18 December 2020
# - Rename sheet "jointUseloginTestAllRoles" -> "jointUseloginTestSubmitter"
# - Change the joint-use submitter's email from ts.allroles@gmail.com to
#   ts.submitter@gmail.com and turn it into a mailto hyperlink (A2),
#   matching the style/column width/selection nudges Excel makes along
#   the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "jointUseloginTestSubmitter"

# Turn A2 into a mailto: hyperlink for the new submitter address. Passing an
# explicit TextToDisplay first lets Excel record the "mailto:..." display
# text on the hyperlink; we then restore the cell text to the plain email
# address (matching what the workbook actually shows).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:ts.submitter@gmail.com", [Type]::Missing, [Type]::Missing, "mailto:ts.submitter@gmail.com") | Out-Null
$ws.Range("A2").Value = "ts.submitter@gmail.com"
$ws.Range("A2").Style = "Hyperlink"

# Column A needs to grow to fit the new (longer) e-mail address.
$ws.Columns("A:A").AutoFit() | Out-Null

# The author's cursor ended up on J13 when they saved.
$ws.Range("J13").Select() | Out-Null
